$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.937.72'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.66%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.906.42'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.52%  '
$ws.Range('E4').Value = '  -0.43%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.74%  '
$ws.Range('E6').Value = '  -0.44%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4837'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.53%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3800'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07374'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9322'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.77'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.54%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07752'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.50%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.945.23'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.491'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.636'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.80'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.85%  '
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008862'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.005'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.997.47'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.99%  '
$ws.Range('E21').Value = '  +0.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.159'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.178.41'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.16%  '
$ws.Range('E24').Value = '  +2.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.917'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.47'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.138'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '117.43'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.960'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08952'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.84%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.230'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.02%  '
$ws.Range('E33').Value = '  +4.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7670'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.99%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.661'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02053'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.537'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.098'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.91%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05283'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5490'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.003'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('E42').Value = '  -0.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1529'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.476'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '109.87'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.63'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4817'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.005'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.655'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.03'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06083'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.14%  '
